$wb = $excel.ActiveWorkbook

# --- "Version" sheet: bump the date stamp ---
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("B2").Value = 44545
$wsVersion.Range("B2").Select()

# --- "Columns" sheet: drop the reference_allele and weight_type rows ---
$wsColumns = $wb.Worksheets.Item("Columns")

# reference_allele / Reference Allele row (row 6)
$wsColumns.Range("A6:C6").Delete()

# weight_type / Type of Weight row (was row 9, now row 8 after the delete above)
$wsColumns.Range("A8:C8").Delete()

$wsColumns.Range("B24").Select()
